$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '20.539.01'
$ws.Range("E2").Value = '  +1.56%  '

# Row 3
$ws.Range("D3").Value = '1.475.27'
$ws.Range("E3").Value = '  +3.00%  '

# Row 4
$ws.Range("E4").Value = '  +0.52%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9690'
$ws.Range("E5").Value = '  -2.61%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '276.97'
$ws.Range("E6").Value = '  -0.25%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3651'
$ws.Range("E7").Value = '  -1.56%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3057'
$ws.Range("E8").Value = '  -3.39%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '40.59'
$ws.Range("E9").Value = '  +0.43%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.060'
$ws.Range("E10").Value = '  -0.19%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06630'
$ws.Range("E11").Value = '  +0.50%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9996'
$ws.Range("E12").Value = '  +0.23%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.479'
$ws.Range("E13").Value = '  -1.68%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.17'
$ws.Range("E14").Value = '  -0.66%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.182'

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001032'
$ws.Range("E16").Value = '  -0.22%  '

# Row 17
$ws.Range("D17").Value = '1.477.38'
$ws.Range("E17").Value = '  +3.13%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9686'
$ws.Range("E18").Value = '  -2.59%  '

# Row 19
$ws.Range("E19").Value = '  +2.71%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.36'
$ws.Range("E20").Value = '  -3.51%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.454'
$ws.Range("E21").Value = '  -3.30%  '

# Row 22
$ws.Range("E22").Value = '  -2.33%  '

# Row 23
$ws.Range("E23").Value = '  -1.04%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.249'
$ws.Range("E24").Value = '  +0.47%  '

# Row 25
$ws.Range("D25").Value = '20.573.28'
$ws.Range("E25").Value = '  +1.64%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.03'
$ws.Range("E26").Value = '  +4.05%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.137'
$ws.Range("E27").Value = '  -8.06%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.23'
$ws.Range("E28").Value = '  -1.48%  '

# Row 29
$ws.Range("D29").Value = '1.632.75'
$ws.Range("E29").Value = '  +2.46%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '113.95'
$ws.Range("E30").Value = '  +1.84%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.933'
$ws.Range("E31").Value = '  -1.01%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8142'
$ws.Range("E32").Value = '  -3.96%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.970'
$ws.Range("E33").Value = '  -6.88%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07917'
$ws.Range("E34").Value = '  +1.31%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.542'
$ws.Range("E35").Value = '  +2.91%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.214'
$ws.Range("E36").Value = '  +8.81%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05820'
$ws.Range("E37").Value = '  -2.06%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.724'
$ws.Range("E38").Value = '  -4.42%  '

# Row 39
$ws.Range("B39").Value = 'Frax'
$ws.Range("C39").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9677'
$ws.Range("E39").Value = '  -2.72%  '

# Row 40
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.676'
$ws.Range("E40").Value = '  -2.30%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02038'
$ws.Range("E41").Value = '  -1.63%  '

# Row 42
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.46'
$ws.Range("E42").Value = '  -3.35%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1882'
$ws.Range("E43").Value = '  -0.55%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5296'
$ws.Range("E44").Value = '  -1.76%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.507'
$ws.Range("E45").Value = '  -1.65%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.15'
$ws.Range("E46").Value = '  -2.01%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '118.06'
$ws.Range("E47").Value = '  -1.19%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5201'
$ws.Range("E48").Value = '  -1.94%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.795'
$ws.Range("E49").Value = '  -0.70%  '

# Row 50
$ws.Range("E50").Value = '  +2.54%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9928'
$ws.Range("E51").Value = '  -0.48%  '
